$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.914
$ws.Range("A9").Value = -20.912
$ws.Range("E11").Value = 12.914
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.738
